$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col18a1"
$ws.Range("C2").Value = "Gpc4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 15.79399266666667
$ws.Range("H2").Value = 47.381978
$ws.Range("I2").Value = 0.2968109173698557
$ws.Range("J2").Value = 0.2968109173698557
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.992380999999999
$ws.Range("N2").Value = 20.977143
$ws.Range("O2").Value = 0.1454502834669897
$ws.Range("P2").Value = 0.1454502834669897
$ws.Range("Q2").Value = 110.4376142365393
$ws.Range("R2").Value = 993.9385281288539
$ws.Range("S2").Value = 0.04317123206754277
$ws.Range("T2").Value = 0.04317123206754276

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col18a1"
$ws.Range("C3").Value = "Gpc4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.79399266666667
$ws.Range("H3").Value = 47.381978
$ws.Range("I3").Value = 0.2968109173698557
$ws.Range("J3").Value = 0.2968109173698557
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 30.83466
$ws.Range("N3").Value = 92.50398
$ws.Range("O3").Value = 0.641399551541635
$ws.Range("P3").Value = 0.641399551541635
$ws.Range("Q3").Value = 487.0023939191601
$ws.Range("R3").Value = 4383.021545272441
$ws.Range("S3").Value = 0.1903743892936867
$ws.Range("T3").Value = 0.1903743892936867

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col18a1"
$ws.Range("C4").Value = "Gpc4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.79399266666667
$ws.Range("H4").Value = 47.381978
$ws.Range("I4").Value = 0.2968109173698557
$ws.Range("J4").Value = 0.2968109173698557
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.246987
$ws.Range("N4").Value = 30.740961
$ws.Range("O4").Value = 0.2131501649913754
$ws.Range("P4").Value = 0.2131501649913754
$ws.Range("Q4").Value = 161.8408375334287
$ws.Range("R4").Value = 1456.567537800858
$ws.Range("S4").Value = 0.06326529600862622
$ws.Range("T4").Value = 0.06326529600862621

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col18a1"
$ws.Range("C5").Value = "Gpc4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.86954866666666
$ws.Range("H5").Value = 74.608646
$ws.Range("I5").Value = 0.4673646309781075
$ws.Range("J5").Value = 0.4673646309781075
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.992380999999999
$ws.Range("N5").Value = 20.977143
$ws.Range("O5").Value = 0.1454502834669897
$ws.Range("P5").Value = 0.1454502834669897
$ws.Range("Q5").Value = 173.8973595753753
$ws.Range("R5").Value = 1565.076236178378
$ws.Range("S5").Value = 0.06797831805821077
$ws.Range("T5").Value = 0.06797831805821077

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col18a1"
$ws.Range("C6").Value = "Gpc4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.86954866666666
$ws.Range("H6").Value = 74.608646
$ws.Range("I6").Value = 0.4673646309781075
$ws.Range("J6").Value = 0.4673646309781075
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 30.83466
$ws.Range("N6").Value = 92.50398
$ws.Range("O6").Value = 0.641399551541635
$ws.Range("P6").Value = 0.641399551541635
$ws.Range("Q6").Value = 766.8440774901198
$ws.Range("R6").Value = 6901.596697411079
$ws.Range("S6").Value = 0.2997674647157799
$ws.Range("T6").Value = 0.2997674647157799

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col18a1"
$ws.Range("C7").Value = "Gpc4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.86954866666666
$ws.Range("H7").Value = 74.608646
$ws.Range("I7").Value = 0.4673646309781075
$ws.Range("J7").Value = 0.4673646309781075
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.246987
$ws.Range("N7").Value = 30.740961
$ws.Range("O7").Value = 0.2131501649913754
$ws.Range("P7").Value = 0.2131501649913754
$ws.Range("Q7").Value = 254.8379418832006
$ws.Range("R7").Value = 2293.541476948806
$ws.Range("S7").Value = 0.09961884820411689
$ws.Range("T7").Value = 0.09961884820411689

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col18a1"
$ws.Range("C8").Value = "Gpc4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.54876233333333
$ws.Range("H8").Value = 37.646287
$ws.Range("I8").Value = 0.2358244516520368
$ws.Range("J8").Value = 0.2358244516520368
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.992380999999999
$ws.Range("N8").Value = 20.977143
$ws.Range("O8").Value = 0.1454502834669897
$ws.Range("P8").Value = 0.1454502834669897
$ws.Range("Q8").Value = 87.74572731311567
$ws.Range("R8").Value = 789.7115458180409
$ws.Range("S8").Value = 0.03430073334123616
$ws.Range("T8").Value = 0.03430073334123616

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col18a1"
$ws.Range("C9").Value = "Gpc4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.54876233333333
$ws.Range("H9").Value = 37.646287
$ws.Range("I9").Value = 0.2358244516520368
$ws.Range("J9").Value = 0.2358244516520368
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 30.83466
$ws.Range("N9").Value = 92.50398
$ws.Range("O9").Value = 0.641399551541635
$ws.Range("P9").Value = 0.641399551541635
$ws.Range("Q9").Value = 386.93681996914
$ws.Range("R9").Value = 3482.43137972226
$ws.Range("S9").Value = 0.1512576975321684
$ws.Range("T9").Value = 0.1512576975321684

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col18a1"
$ws.Range("C10").Value = "Gpc4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.54876233333333
$ws.Range("H10").Value = 37.646287
$ws.Range("I10").Value = 0.2358244516520368
$ws.Range("J10").Value = 0.2358244516520368
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.246987
$ws.Range("N10").Value = 30.740961
$ws.Range("O10").Value = 0.2131501649913754
$ws.Range("P10").Value = 0.2131501649913754
$ws.Range("Q10").Value = 128.5870044957563
$ws.Range("R10").Value = 1157.283040461807
$ws.Range("S10").Value = 0.05026602077863228
$ws.Range("T10").Value = 0.05026602077863227
